$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to stay text so values like "22.30" keep their
# trailing zero instead of being auto-coerced to a number by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.403.46"
$ws.Range("E2").Value = "  -2.98%  "
$ws.Range("D3").Value = "2.245.94"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "234.57"
$ws.Range("D6").Value = "0.633"
$ws.Range("E6").Value = "  -4.43%  "
$ws.Range("D7").Value = "69.73"
$ws.Range("E7").Value = "  -2.96%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.561"
$ws.Range("E9").Value = "  -4.74%  "
$ws.Range("D10").Value = "0.0997"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("D11").Value = "58.82"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D12").Value = "36.54"
$ws.Range("E12").Value = "  +12.89%  "
$ws.Range("E13").Value = "  -1.54%  "
$ws.Range("D14").Value = "6.74"
$ws.Range("E14").Value = "  -5.10%  "
$ws.Range("D15").Value = "2.581.63"
$ws.Range("E15").Value = "  -3.71%  "
$ws.Range("D16").Value = "15.13"
$ws.Range("E16").Value = "  -5.53%  "
$ws.Range("D17").Value = "0.855"
$ws.Range("E17").Value = "  -4.20%  "
$ws.Range("D18").Value = "2.251.13"
$ws.Range("E18").Value = "  -3.67%  "
$ws.Range("D19").Value = "42.272.10"
$ws.Range("E19").Value = "  -3.07%  "
$ws.Range("D20").Value = "0.0₃0979"
$ws.Range("E20").Value = "  -2.37%  "
$ws.Range("D21").Value = "6.26"
$ws.Range("E21").Value = "  -4.62%  "
$ws.Range("D22").Value = "73.47"
$ws.Range("D23").Value = "236.41"
$ws.Range("E23").Value = "  -4.61%  "
$ws.Range("D24").Value = "1.97"
$ws.Range("E24").Value = "  +3.66%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "3.69"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("D27").Value = "2.41"
$ws.Range("E27").Value = "  -3.35%  "
$ws.Range("D28").Value = "10.02"
$ws.Range("E28").Value = "  -2.76%  "
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  -2.02%  "
$ws.Range("D30").Value = "170.35"
$ws.Range("E30").Value = "  -2.94%  "
$ws.Range("D31").Value = "20.58"
$ws.Range("E31").Value = "  -6.89%  "
$ws.Range("E32").Value = "  -3.67%  "
$ws.Range("E33").Value = "  -5.29%  "
$ws.Range("D34").Value = "0.0729"
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("D35").Value = "5.35"
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").Value = "4.72"
$ws.Range("E36").Value = "  -6.41%  "
$ws.Range("D37").Value = "3.66"
$ws.Range("E37").Value = "  -1.53%  "
$ws.Range("D38").Value = "22.30"
$ws.Range("E38").Value = "  +20.16%  "
$ws.Range("D39").Value = "0.0278"
$ws.Range("E39").Value = "  +3.22%  "
$ws.Range("D40").Value = "2.29"
$ws.Range("E40").Value = "  -3.00%  "
$ws.Range("D41").Value = "5.98"
$ws.Range("E41").Value = "  -6.17%  "
$ws.Range("D42").Value = "65.45"
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("D43").Value = "9.23"
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("D44").Value = "4.98"
$ws.Range("E44").Value = "  -12.66%  "
$ws.Range("E45").Value = "  -2.56%  "
$ws.Range("E46").Value = "  -1.87%  "
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").Value = "4.54"
$ws.Range("E48").Value = "  +11.82%  "
$ws.Range("D49").Value = "10.25"
$ws.Range("E49").Value = "  +10.23%  "
$ws.Range("E50").Value = "  -2.67%  "
$ws.Range("D51").Value = "2.35"
$ws.Range("E51").Value = "  -2.69%  "

# Restore column D formatting/style so no stray format is left on the cells
$ws.Range("D2:D51").NumberFormat = "General"
$ws.Range("D2:D51").Style = "Normal"
